# Update Ntn1-Unc5d LR-pair sheet with recomputed TPM-based NATMI output.
# - rows 2..5 get recalculated values (columns A, D, G..T); B/C/E/F unchanged
# - rows 6..9 (the former MuSCs / Resolving-Mac duplicate rows) are removed

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Sending cluster" (A) / "Target cluster" (D) labels for rows 2-5
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 4).Value = "MuSCs"

$ws.Cells.Item(3, 1).Value = "FAPs"
$ws.Cells.Item(3, 4).Value = "ECs"

$ws.Cells.Item(4, 1).Value = "MuSCs"
$ws.Cells.Item(4, 4).Value = "MuSCs"

$ws.Cells.Item(5, 1).Value = "Resolving-Mac"
$ws.Cells.Item(5, 4).Value = "MuSCs"

# Recomputed numeric columns G..T (cols 7..20) for rows 2-5 (B, C, E, F unchanged)
$data = @(
    @(2.629231666666667, 7.887695, 0.1414315557047068, 0.1414315557047067, 2, 0.6666666666666666, 0.05794366666666667, 0.173831, 1, 1, 0.1523473232827778, 1.371125909545, 0.1414315557047068, 0.1414315557047067),
    @(11.42765333333333, 34.28296, 0.6147160060020365, 0.6147160060020365, 2, 0.6666666666666666, 0.05794366666666667, 0.173831, 1, 1, 0.662160135528889, 5.959441219760001, 0.6147160060020365, 0.6147160060020365),
    @(4.24731, 12.74193, 0.2284711798035388, 0.2284711798035388, 2, 0.6666666666666666, 0.05794366666666667, 0.173831, 1, 1, 0.24610471487, 2.21494243383, 0.2284711798035388, 0.2284711798035388),
    @(0.2859396666666667, 0.8578190000000001, 0.01538125848971795, 0.01538125848971795, 2, 0.6666666666666666, 0.05794366666666667, 0.173831, 1, 1, 0.01656839273211112, 0.149115534589, 0.01538125848971795, 0.01538125848971795)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $values = $data[$i]
    for ($j = 0; $j -lt $values.Length; $j++) {
        $col = 7 + $j
        $ws.Cells.Item($row, $col).Value = $values[$j]
    }
}

# Drop the four rows that no longer exist in the refreshed export (delete
# from the bottom up so row numbers of the still-pending rows stay valid)
$ws.Rows.Item(9).Delete()
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(6).Delete()
